$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.388317276603885
$ws.Range("C2").Value = 0.04981286588697742
$ws.Range("D2").Value = 0.1577979911730818
$ws.Range("E2").Value = 0.07019828564043351
$ws.Range("F2").Value = 2.625751773257264
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.9305988667144902
$ws.Range("L2").Value = 0.1963251175074916
$ws.Range("M2").Value = 0.298867986704245
$ws.Range("N2").Value = 3.740235092271789
$ws.Range("B3").Value = 1.347845959518565
$ws.Range("C3").Value = 0.04514965112339553
$ws.Range("D3").Value = 0.1581448829324259
$ws.Range("E3").Value = 0.07033141536888365
$ws.Range("F3").Value = 2.602581907047337
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.8883912467549351
$ws.Range("L3").Value = 0.1940149243992977
$ws.Range("M3").Value = 0.2919269069122592
$ws.Range("N3").Value = 3.7423052222113
$ws.Range("B4").Value = 1.323817732214849
$ws.Range("C4").Value = 0.04226370416380121
$ws.Range("D4").Value = 0.1583614021338491
$ws.Range("E4").Value = 0.07043031795491217
$ws.Range("F4").Value = 2.589425679108842
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.8630588559736907
$ws.Range("L4").Value = 0.1926930427625635
$ws.Range("M4").Value = 0.2878346556715989
$ws.Range("N4").Value = 3.744342814676344
$ws.Range("B5").Value = 1.314232525788867
$ws.Range("C5").Value = 0.04108186791624746
$ws.Range("D5").Value = 0.1584505233941647
$ws.Range("E5").Value = 0.07047494563243806
$ws.Range("F5").Value = 2.584333303512111
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.852882165177931
$ws.Range("L5").Value = 0.1921786667685197
$ws.Range("M5").Value = 0.2862096803863317
$ws.Range("N5").Value = 3.745365747169942
$ws.Range("B6").Value = 1.312653380257615
$ws.Range("C6").Value = 0.0408852729739948
$ws.Range("D6").Value = 0.1584653755781034
$ws.Range("E6").Value = 0.07048261745194839
$ws.Range("F6").Value = 2.583503953079287
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.8512011738118019
$ws.Range("L6").Value = 0.1920947232103245
$ws.Range("M6").Value = 0.2859424305853224
$ws.Range("N6").Value = 3.745547232352379
$ws.Range("B7").Value = 1.323687626740792
$ws.Range("C7").Value = 0.04224778905032167
$ws.Range("D7").Value = 0.1583626004583856
$ws.Range("E7").Value = 0.07043090229925397
$ws.Range("F7").Value = 2.589355913093144
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.8629210167962924
$ws.Range("L7").Value = 0.1926860072924796
$ws.Range("M7").Value = 0.2878125679921411
$ws.Range("N7").Value = 3.744355830705871
$ws.Range("B8").Value = 1.37419225928474
$ws.Range("C8").Value = 0.04820965574594993
$ws.Range("D8").Value = 0.1579168683519878
$ws.Range("E8").Value = 0.07024063242690026
$ws.Range("F8").Value = 2.617540552107741
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.9159244858467162
$ws.Range("L8").Value = 0.1955085235183063
$ws.Range("M8").Value = 0.2964395004088871
$ws.Range("N8").Value = 3.740789685325097
$ws.Range("B9").Value = 1.479757675600808
$ws.Range("C9").Value = 0.05972461984394783
$ws.Range("D9").Value = 0.1570706872282477
$ws.Range("E9").Value = 0.0700033073951456
$ws.Range("F9").Value = 2.681316881054997
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 1.024507976927111
$ws.Range("L9").Value = 0.2018097738820828
$ws.Range("M9").Value = 0.3147038105791324
$ws.Range("N9").Value = 3.739887086963947
$ws.Range("B10").Value = 1.56131840834729
$ws.Range("C10").Value = 0.06808343973087005
$ws.Range("D10").Value = 0.1564658558218222
$ws.Range("E10").Value = 0.06991127517287055
$ws.Range("F10").Value = 2.733388172271361
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 1.1071473135205
$ws.Range("L10").Value = 0.2069072007137578
$ws.Range("M10").Value = 0.3289473801535863
$ws.Range("N10").Value = 3.742951321020215
$ws.Range("B11").Value = 1.599297271945886
$ws.Range("C11").Value = 0.07186548398262005
$ws.Range("D11").Value = 0.1561943312605463
$ws.Range("E11").Value = 0.06988719263793364
$ws.Range("F11").Value = 2.758215998717304
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 1.145371365710048
$ws.Range("L11").Value = 0.2093279736353679
$ws.Range("M11").Value = 0.3356071743520346
$ws.Range("N11").Value = 3.745158021733459
$ws.Range("B12").Value = 1.613805177091365
$ws.Range("C12").Value = 0.07329480452794712
$ws.Range("D12").Value = 0.1560920308510987
$ws.Range("E12").Value = 0.06988062239082993
$ws.Range("F12").Value = 2.767782036644519
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 1.159936972332588
$ws.Range("L12").Value = 0.2102593172544545
$ws.Range("M12").Value = 0.3381550330237459
$ws.Range("N12").Value = 3.746110755368363
$ws.Range("B13").Value = 1.610675028672574
$ws.Range("C13").Value = 0.07298710059343705
$ws.Range("D13").Value = 0.1561140399431791
$ws.Range("E13").Value = 0.06988192415656869
$ws.Range("F13").Value = 2.765714508011769
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 1.156795957580499
$ws.Range("L13").Value = 0.2100580843192006
$ws.Range("M13").Value = 0.3376051522638548
$ws.Range("N13").Value = 3.745900355197108
$ws.Range("B14").Value = 1.600488317540396
$ws.Range("C14").Value = 0.07198313185689642
$ws.Range("D14").Value = 0.156185904521787
$ws.Range("E14").Value = 0.06988660105349709
$ws.Range("F14").Value = 2.758999708277514
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 1.146567861948768
$ws.Range("L14").Value = 0.2094043022797081
$ws.Range("M14").Value = 0.3358162683029917
$ws.Range("N14").Value = 3.745234055436697
$ws.Range("B15").Value = 1.594265092195258
$ws.Range("C15").Value = 0.07136780224887218
$ws.Range("D15").Value = 0.1562299913804948
$ws.Range("E15").Value = 0.06988979754599711
$ws.Range("F15").Value = 2.754908102541776
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 1.140314714856657
$ws.Range("L15").Value = 0.2090057494784361
$ws.Range("M15").Value = 0.3347239035558758
$ws.Range("N15").Value = 3.744841184867681
$ws.Range("B16").Value = 1.55885404540652
$ws.Range("C16").Value = 0.06783587003405955
$ws.Range("D16").Value = 0.1564836733584025
$ws.Range("E16").Value = 0.06991320618941721
$ws.Range("F16").Value = 2.731788587368897
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 1.104662000445103
$ws.Range("L16").Value = 0.2067510478909185
$ws.Range("M16").Value = 0.3285157762293807
$ws.Range("N16").Value = 3.742823482351554
$ws.Range("B17").Value = 1.537355068878753
$ws.Range("C17").Value = 0.0656639691974874
$ws.Range("D17").Value = 0.1566402252239634
$ws.Range("E17").Value = 0.06993211570711999
$ws.Range("F17").Value = 2.717897792047225
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 1.082952005911636
$ws.Range("L17").Value = 0.205393962020409
$ws.Range("M17").Value = 0.3247534821262477
$ws.Range("N17").Value = 3.741794017593818
$ws.Range("B18").Value = 1.525071918107642
$ws.Range("C18").Value = 0.06441282201710408
$ws.Range("D18").Value = 0.1567306100859582
$ws.Range("E18").Value = 0.06994466628777829
$ws.Range("F18").Value = 2.710015489341146
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 1.070524368944064
$ws.Range("L18").Value = 0.2046229964757202
$ws.Range("M18").Value = 0.3226064822985393
$ws.Range("N18").Value = 3.741278384042715
$ws.Range("B19").Value = 1.520927219525731
$ws.Range("C19").Value = 0.06398887165359213
$ws.Range("D19").Value = 0.1567612713068121
$ws.Range("E19").Value = 0.06994920350929057
$ws.Range("F19").Value = 2.707365099557421
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 1.066326775262553
$ws.Range("L19").Value = 0.2043636088725123
$ws.Range("M19").Value = 0.321882459794935
$ws.Range("N19").Value = 3.741116929230174
$ws.Range("B20").Value = 1.539635133361742
$ws.Range("C20").Value = 0.0658953706361558
$ws.Range("D20").Value = 0.1566235247848002
$ws.Range("E20").Value = 0.06992992952809551
$ws.Range("F20").Value = 2.71936538207143
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 1.085256924740662
$ws.Range("L20").Value = 0.2055374332065867
$ws.Range("M20").Value = 0.3251522281363037
$ws.Range("N20").Value = 3.7418956881219
$ws.Range("B21").Value = 1.60347697748125
$ws.Range("C21").Value = 0.07227809872792079
$ws.Range("D21").Value = 0.1561647820382284
$ws.Range("E21").Value = 0.06988515821155872
$ws.Range("F21").Value = 2.760967546818847
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 1.149569630163967
$ws.Range("L21").Value = 0.2095959363886664
$ws.Range("M21").Value = 0.3363410028263445
$ws.Range("N21").Value = 3.745426583808268
$ws.Range("B22").Value = 1.645936567904187
$ws.Range("C22").Value = 0.07643296182878601
$ws.Range("D22").Value = 0.1558679984995237
$ws.Range("E22").Value = 0.06987075313221425
$ws.Range("F22").Value = 2.78911465389325
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 1.192132167361194
$ws.Range("L22").Value = 0.2123337879644396
$ws.Range("M22").Value = 0.3438047122480157
$ws.Range("N22").Value = 3.748416909460786
$ws.Range("B23").Value = 1.623207775367746
$ws.Range("C23").Value = 0.0742169274207356
$ws.Range("D23").Value = 0.1560261200393542
$ws.Range("E23").Value = 0.06987708471126197
$ws.Range("F23").Value = 2.774004283645354
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 1.16936711812923
$ws.Range("L23").Value = 0.2108647346637156
$ws.Range("M23").Value = 0.3398073520782887
$ws.Range("N23").Value = 3.746758374091314
$ws.Range("B24").Value = 1.538604076693332
$ws.Range("C24").Value = 0.06579076181627386
$ws.Range("D24").Value = 0.1566310738647845
$ws.Range("E24").Value = 0.06993091266840068
$ws.Range("F24").Value = 2.718701561623448
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 1.084214703604403
$ws.Range("L24").Value = 0.2054725410985299
$ws.Range("M24").Value = 0.3249719052499245
$ws.Range("N24").Value = 3.741849485464172
$ws.Range("B25").Value = 1.450498390351584
$ws.Range("C25").Value = 0.05662774336107645
$ws.Range("D25").Value = 0.157296631002211
$ws.Range("E25").Value = 0.07005301959852339
$ws.Range("F25").Value = 2.663150046952751
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.9946328228310506
$ws.Range("L25").Value = 0.2000230065991246
$ws.Range("M25").Value = 0.3096182602906197
$ws.Range("N25").Value = 3.739477753211426